# Update Name of Algo
# Applies corrected imputed values to the result data produced by the
# RandomForest algorithm (terrestrial_mammals / combination_2_ABCDE / ACE /
# seed 2). Only specific data cells change value; everything else
# (headers, formatting, other cells) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -22.09570000000001
$ws.Range("E3").Value = 16.41409999999999
$ws.Range("A21").Value = -19.92659999999999
$ws.Range("A23").Value = -20.55419999999998
$ws.Range("E24").Value = 16.46540000000001
$ws.Range("A25").Value = -21.58369999999999
$ws.Range("C27").Value = -12.59299999999999
$ws.Range("C31").Value = -13.0025
$ws.Range("C39").Value = -12.5534
$ws.Range("C48").Value = -11.20939999999999
$ws.Range("C51").Value = -11.5896
$ws.Range("C52").Value = -10.9459
$ws.Range("A53").Value = -21.8694
$ws.Range("C55").Value = -13.4202
$ws.Range("C56").Value = -12.53029999999999
$ws.Range("A57").Value = -22.01079999999999
$ws.Range("C57").Value = -12.94299999999999
$ws.Range("E57").Value = 16.55159999999999
$ws.Range("A59").Value = -21.981
$ws.Range("E61").Value = 16.4627
$ws.Range("A69").Value = -21.5935
$ws.Range("E70").Value = 17.39460000000001
$ws.Range("C73").Value = -12.8124
$ws.Range("A79").Value = -20.63840000000002
$ws.Range("A83").Value = -21.8757
$ws.Range("E86").Value = 16.70170000000001
$ws.Range("C89").Value = -10.2664
$ws.Range("C90").Value = -12.6022
$ws.Range("A93").Value = -21.57449999999999
$ws.Range("E98").Value = 15.9499
$ws.Range("E100").Value = 16.44190000000001
$ws.Range("E102").Value = 16.78569999999999
